$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "jak push..." text from C17 entirely and restore default row height
$ws.Range("C17").Clear()
$ws.Range("B17:C17").EntireRow.AutoFit()

# Add new row 20 - copy style from row 17's date cell, then set values
$ws.Range("B17").Copy($ws.Range("B20"))
$ws.Range("B20").Value = 43336
$ws.Range("C20").Value = "sec5lec37 komentarze"

# Update the view state (scroll position + selection) to match final editing state
$ws.Range("C21").Select()
